# DhalisMenu.xlsx - add a new menu item row to Sheet1
#
# A new row is inserted at row 9 ("7 Inch Pizza Veggi Panner( Onion and
# Capcium and corn)", Full price 120, no Half price) which pushes all the
# following rows (previously 9-52, now 10-53) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 9, shifting rows 9..52 down to 10..53.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 (note: column B / "Half" is intentionally left blank).
$ws.Range("A9").Value = "7 Inch Pizza Veggi Panner( Onion and Capcium and corn)"
$ws.Range("C9").Value = 120
$ws.Range("D9").Value = "7 Inch Pizza Veggi Panner( Onion and Capcium and corn).png"

# Update the selected cell to D10 (matches the saved view state in the workbook).
$null = $ws.Range("D10").Select()
